# Week 4 logs: copy over Jesse's local activity/task data into the
# Task Summary Sheet and the Activity Log Summary Sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# TASK SUMMARY SHEET
# ---------------------------------------------------------------------------
$taskSummary = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Header: name + week number
$taskSummary.Range("C1").Value = "Jesse Hare"
$taskSummary.Range("E1").Value = 4

# Stage / Task rows
$taskSummary.Range("A3").Value = "Project Build"
$taskSummary.Range("B3").Value = "Optimising searching/sorting functions"
$taskSummary.Range("C3").Value = 4
$taskSummary.Range("D3").Value = 5

$taskSummary.Range("A4").Value = "Project Build"
$taskSummary.Range("B4").Value = "Implementing new requirements"
$taskSummary.Range("C4").Value = 3
$taskSummary.Range("D4").Value = 7

$taskSummary.Range("A5").Value = "Project Build"
$taskSummary.Range("B5").Value = "Unit testing of program modules"
$taskSummary.Range("C5").Value = 4
$taskSummary.Range("D5").Value = 3

$taskSummary.Range("A6").Value = "Project Build"
$taskSummary.Range("B6").Value = "Testing of GUI on different types of devices"
$taskSummary.Range("C6").Value = 2
$taskSummary.Range("D6").Value = 3

$taskSummary.Range("A7").Value = "Project Build"
$taskSummary.Range("B7").Value = "Restructuring of code into a more logical layout"
$taskSummary.Range("C7").Value = 2
$taskSummary.Range("D7").Value = 2

# ---------------------------------------------------------------------------
# ACTIVITY LOG SUMMARY SHEET
# ---------------------------------------------------------------------------
$logSummary = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")

$logSummary.Range("D1").Value = "Jesse Hare"

$logSummary.Range("A4").Value = "Project Build"
$logSummary.Range("B4").Value = 12
$logSummary.Range("C4").Value = 8

$wb.Application.CalculateFull()
